# Initial version with clustering, aggregated and extreme periods.
# Updates a few input values on the "thermal", "cap_factors" and "demand"
# sheets, and leaves the workbook with the "thermal" tab active/selected
# (cell B3), matching the author's last editing session.

$wb = $excel.ActiveWorkbook

# --- Data edits -----------------------------------------------------------

# thermal!B2 : max_cap for t1, 1000 -> 500
$thermal = $wb.Worksheets.Item("thermal")
$thermal.Range("B2").Value = 500

# cap_factors!C2 : cap_factor for w1, 0.7 -> 0.684842
$capFactors = $wb.Worksheets.Item("cap_factors")
$capFactors.Range("C2").Value = 0.684842

# demand!B2 : demand for period 1, 500 -> 477.227212
$demand = $wb.Worksheets.Item("demand")
$demand.Range("B2").Value = 477.227212

# --- View / selection state -------------------------------------------

# Leave the "thermal" worksheet as the active tab with B3 selected.
$thermal.Activate()
$thermal.Range("B3").Select()

# Reposition/resize the Excel window to match the author's session.
$aw = $excel.ActiveWindow
$aw.Left = -26988
$aw.Top = 516
$aw.Width = 17280
$aw.Height = 9852
